$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.286.45"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.497.61"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.75"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.25"
$ws.Range("E6").Value = "  +3.02%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +1.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.21"
$ws.Range("E10").Value = "  +0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.388"
$ws.Range("E11").Value = "  +2.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.088.29"
$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("E14").Value = "  +2.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.494.50"
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.87"
$ws.Range("E16").Value = "  -5.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.321.28"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.80"
$ws.Range("E19").Value = "  +3.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.80"
$ws.Range("E20").Value = "  -3.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.53"
$ws.Range("E21").Value = "  +1.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.569"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.636.21"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.14"
$ws.Range("E24").Value = "  +1.79%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +1.56%  "

$ws.Range("E27").Value = "  +2.37%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.45"
$ws.Range("E28").Value = "  +0.94%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.29"
$ws.Range("E30").Value = "  +1.68%  "

$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("E32").Value = "  -4.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.516.69"
$ws.Range("E33").Value = "  +0.85%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("E35").Value = "  +4.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.55"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.22"
$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.89"
$ws.Range("E39").Value = "  +0.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.61"
$ws.Range("E40").Value = "  -2.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0783"
$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.92"
$ws.Range("E43").Value = "  -2.77%  "

$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.76"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.43"
$ws.Range("E46").Value = "  +1.58%  "

$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("E48").Value = "  -2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.483.25"
$ws.Range("E49").Value = "  +2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("E51").Value = "  +1.68%  "
